# Update CV/resume job-title text in column B for several date ranges.
# The author renamed two overlapping "Finance participative" roles at La Nef /
# Artisans Angkor to "Chargée de projet Innovation", and renamed the later
# "Attachée Territoriale" role (from 2021-06 onward) to mention
# "Politique médico-sociale locale" instead of "relations partenariales".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$stageNef    = "Je réalise un stage en tant que Chargée de projet Innovation à la Nef"
$stageAngkor = "Je réalise mon stage de fin d'étude en qualité de Chargée de projet Innovation chez Artisans Angkor."
$nefNouveau  = "Je travaille de nouveau à la Nef, cette fois en tant que Chargée de projet Innovation"
$attachee    = "J'exerce la fonction d'Attachée Territoriale Politique médico-sociale locale au Département du Rhône. "

foreach ($r in 47..52) {
    $ws.Range("B$r").Value = $stageNef
}

foreach ($r in 56..61) {
    $ws.Range("B$r").Value = $stageAngkor
}

foreach ($r in 62..71) {
    $ws.Range("B$r").Value = $nefNouveau
}

foreach ($r in 73..83) {
    $ws.Range("B$r").Value = $attachee
}

$ws.Range("C87").Select()
